$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (pushes old row19.. content -- here nothing -- and
# rows 23/24 -- down by one, preserving their formatting/merges).
$ws.Rows("19").Insert()

# Copy row 18 (period 2508) down into the freshly inserted row 19 so it
# picks up the same "inner" row formatting that rows 16/17 use, and make the
# old row 18 revert to that same inner style (matches how the template grows
# the table when a new period is appended).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4104) | Out-Null
$ws.Range("B18:J18").Copy()
$ws.Range("B18:J18").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# New period row: same worker/salary data as row 18, new "Periodo Mora" 2509.
$ws.Range("E19").Value2 = "2509"
$ws.Range("F19").Value2 = 70000
$ws.Range("G19").Value2 = 1750000

# Column E ("Tipo Doc Trabajador") of the whole data block is now centered.
$ws.Range("E16:E19").HorizontalAlignment = -4108

# Totals: one more "mora" period (3 -> 4 rows of 70000) and the period count.
$ws.Range("E11").Value2 = 280000
$ws.Range("F13").Value2 = 3
